$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.594.11"
$ws.Range("E2").Value = "  +4.51%  "

$ws.Range("D3").Value = "3.614.28"
$ws.Range("E3").Value = "  +4.51%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "629.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.48"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +7.30%  "

$ws.Range("D7").Value = "3.612.82"
$ws.Range("E7").Value = "  +4.50%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.50"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +9.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.441"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.61"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +8.09%  "

$ws.Range("D15").Value = "4.228.38"
$ws.Range("E15").Value = "  +4.64%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "69.618.53"
$ws.Range("E16").Value = "  +4.70%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.609.77"
$ws.Range("E17").Value = "  +4.68%  "

$ws.Range("E18").Value = "  +0.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.74"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +6.44%  "

$ws.Range("E20").Value = "  +8.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.21"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +14.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "461.63"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.645"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.89"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000134"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +9.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.74"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +7.50%  "

$ws.Range("D27").Value = "3.757.77"
$ws.Range("E27").Value = "  +4.47%  "

$ws.Range("E28").Value = "  +0.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.29"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +13.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.65"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.09%  "

$ws.Range("E31").Value = "  +11.55%  "

$ws.Range("E32").Value = "  +9.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.55"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.62%  "

$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.53"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.95"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.09%  "

$ws.Range("D37").Value = "3.613.40"
$ws.Range("E37").Value = "  +4.97%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.49"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +8.31%  "

$ws.Range("E39").Value = "  +13.65%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0927"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +8.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "178.80"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.90%  "

$ws.Range("E43").Value = "  +0.21%  "

$ws.Range("E44").Value = "  +4.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.17"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +24.54%  "

$ws.Range("E46").Value = "  +4.38%  "

$ws.Range("E47").Value = "  +14.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.90"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.75"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +11.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.82"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.19%  "

$ws.Range("E51").Value = "  +9.86%  "
